$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Continue Python App Course"
$ws.Range("A8").Value = "Chill Day 😁"
$ws.Range("A9").Value = "1H Py Apps Course"
